$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.048.75'
$ws.Range('E2').Value = '  +0.76%  '

$ws.Range('D3').Value = '2.264.78'
$ws.Range('E3').Value = '  -0.38%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.86'
$ws.Range('E5').Value = '  -0.07%  '

$ws.Range('E6').Value = '  +1.76%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.65'
$ws.Range('E7').Value = '  +4.82%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.448'
$ws.Range('E9').Value = '  +5.69%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.100'
$ws.Range('E10').Value = '  +5.41%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.32'
$ws.Range('E11').Value = '  -1.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '27.34'
$ws.Range('E12').Value = '  +14.69%  '

$ws.Range('E13').Value = '  +1.88%  '

$ws.Range('D14').Value = '2.601.97'
$ws.Range('E14').Value = '  -0.45%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.75'
$ws.Range('E15').Value = '  +0.13%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.09'
$ws.Range('E16').Value = '  +4.99%  '

$ws.Range('E17').Value = '  +3.53%  '

$ws.Range('D18').Value = '2.269.23'
$ws.Range('E18').Value = '  -0.76%  '

$ws.Range('D19').Value = '43.948.73'
$ws.Range('E19').Value = '  +0.58%  '

$ws.Range('E20').Value = '  +7.89%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.76'
$ws.Range('E21').Value = '  +1.08%  '

$ws.Range('E22').Value = '  -1.93%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.02'
$ws.Range('E23').Value = '  +0.11%  '

$ws.Range('E24').Value = '  +0.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('E25').Value = '  -3.52%  '

$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.12'
$ws.Range('E26').Value = '  +2.63%  '

$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.25'
$ws.Range('E27').Value = '  -8.96%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.24'
$ws.Range('E28').Value = '  +21.67%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.55'
$ws.Range('E29').Value = '  +0.66%  '

$ws.Range('E30').Value = '  +0.34%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.96'
$ws.Range('E31').Value = '  +1.84%  '

$ws.Range('E32').Value = '  -4.10%  '

$ws.Range('E33').Value = '  +1.93%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0708'
$ws.Range('E34').Value = '  +7.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.81'
$ws.Range('E35').Value = '  +0.30%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.92'
$ws.Range('E36').Value = '  -3.84%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.82'
$ws.Range('E37').Value = '  +6.17%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.53'
$ws.Range('E38').Value = '  +1.10%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.32'
$ws.Range('E39').Value = '  -3.60%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0260'
$ws.Range('E40').Value = '  +3.96%  '

$ws.Range('E41').Value = '  +0.13%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000222'
$ws.Range('E42').Value = '  -2.67%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0983'
$ws.Range('E43').Value = '  -1.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.52'
$ws.Range('E44').Value = '  +5.28%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.49'
$ws.Range('E45').Value = '  +9.16%  '

$ws.Range('E46').Value = '  -5.75%  '

$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '98.48'
$ws.Range('E48').Value = '  +0.37%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.44'
$ws.Range('E49').Value = '  -1.28%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.39'
$ws.Range('E50').Value = '  +5.38%  '

$ws.Range('D51').Value = '1.448.16'
$ws.Range('E51').Value = '  -1.58%  '
